$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E) - updated rep counts
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - meanEMG/legmaxROM recalculated; C2 & E2 now blank
$ws.Range("B2").Value = 6.5925963874547371
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 4.1416649407512338
$ws.Range("E2").ClearContents()

# Row 3 (STR) - meanEMG/legmaxROM recalculated
$ws.Range("B3").Value = 5.4026204843154222
$ws.Range("C3").Value = 6.3751365426387139
$ws.Range("D3").Value = 4.0442631587009723
$ws.Range("E3").Value = 9.9670937305617571

# Selection now only covers the updated B1:E3 block
$ws.Range("B1:E3").Select() | Out-Null
